# "Generate Report for Handoff" -- the handoff status moves from
# "In Translation" to "Ready for handoff" and the relevant timestamps are
# refreshed to reflect the new handoff-generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status + the overall "Latest HO Xliff
#     Generate Date" column.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-26 08:09:25"

# --- zh-cn detail sheet: Status + Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-26 08:09:13"

# --- de-de detail sheet: Status + Latest Handoff Datetime.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-26 08:09:25"

# The "Status" columns are widened so the new, longer "Ready for handoff"
# text fits (matches Excel auto-resizing the column after the value grew).
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
